$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q1" right before the "总计" sheet,
#    carrying fund-level detail (same shape as the 2020-Q4 / 2021-Q4
#    sheets).
# ------------------------------------------------------------------
$beforeAnchor = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($beforeAnchor)
$q1.Name = "2022-Q1"

# NOTE: inserting a sheet shifts everybody after it, so the handle we
# grabbed for "总计" before the Add() is now stale (it resolves to
# whatever sheet now sits at that old index - i.e. our new sheet).
# Re-resolve "总计" by name AFTER the insert so later writes land on
# the right tab.
$totalSheet = $wb.Worksheets.Item("总计")

# Pull the header + index-column formatting from the existing
# "2021-Q4" sheet (style carries border/bold/center + default style).
# Copy B1:H1 (header) and A2:A4 (index column) separately so we don't
# touch A1, which stays empty/unused just like on the sibling sheets.
$src = $wb.Worksheets.Item("2021-Q4")
$src.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$src.Range("A2:A4").Copy()
$q1.Range("A2:A4").PasteSpecial(-4122) | Out-Null

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund code (B) and the decimal-looking metrics (D:G) are all stored
# as *text* in the source workbook (leading zeros in the fund code
# must be preserved, and the metrics are text too) - force text format
# before writing them so Excel doesn't coerce them into numbers.
$q1.Range("B2:B4").NumberFormat = "@"
$q1.Range("D2:G4").NumberFormat = "@"

# Row 2
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "006199"
$q1.Range("C2").Value = "长盛同锦研究精选混合"
$q1.Range("D2").Value = "1.73"
$q1.Range("E2").Value = "82.48"
$q1.Range("F2").Value = "3.00"
$q1.Range("G2").Value = "0.0519"
$q1.Range("H2").Value = 6

# Row 3
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "001892"
$q1.Range("C3").Value = "长盛新兴成长主题灵活配置混合"
$q1.Range("D3").Value = "1.32"
$q1.Range("E3").Value = "82.10"
$q1.Range("F3").Value = "3.03"
$q1.Range("G3").Value = "0.0400"
$q1.Range("H3").Value = 7

# Row 4
$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "002085"
$q1.Range("C4").Value = "长盛互联网+主题灵活配置混合"
$q1.Range("D4").Value = "0.84"
$q1.Range("E4").Value = "83.97"
$q1.Range("F4").Value = "3.00"
$q1.Range("G4").Value = "0.0252"
$q1.Range("H4").Value = 7

$q1.Range("A1").Select()

# ------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: add the new 2022-Q1 row on top
#    of the existing history (2021-Q4, 2020-Q4).
# ------------------------------------------------------------------
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.12

$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 3
$totalSheet.Range("D3").Value = 0.29

$totalSheet.Range("B4").Value = "2020-Q4"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.01

# Row 4 is brand new (the sheet only had 2 data rows before), so carry
# over the index-column formatting (bold/border/center) from row 3
# before setting its value.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122) | Out-Null
$totalSheet.Range("A4").Value = 2

Write-Output "ok"
